$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.823.13"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.641.00"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'607.54"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "'147.14"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("E10").Value = "  +5.23%  "
$ws.Range("D11").Value = "'5.60"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'27.46"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "3.114.89"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "63.662.37"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "2.631.98"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "'11.78"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("E19").Value = "  +3.42%  "
$ws.Range("D20").Value = "'347.23"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "'66.29"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  +8.05%  "
$ws.Range("E26").Value = "  +2.85%  "
$ws.Range("D27").Value = "'9.22"
$ws.Range("E27").Value = "  +5.93%  "
$ws.Range("D28").Value = "'565.63"
$ws.Range("E28").Value = "  +4.00%  "
$ws.Range("D29").Value = "'8.11"
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.160"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("E33").Value = "  +4.94%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'5.29"
$ws.Range("E35").Value = "  +2.30%  "
$ws.Range("D36").Value = "'169.55"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.404"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").Value = "'1.95"
$ws.Range("E39").Value = "  +4.76%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'165.29"
$ws.Range("E42").Value = "  -6.26%  "
$ws.Range("D43").Value = "'40.13"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").Value = "'21.97"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "'0.627"
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("E48").Value = "  +13.66%  "
$ws.Range("D51").Value = "'18.84"
$ws.Range("E51").Value = "  -0.65%  "
